# Fix inconsistent casing in the "key" column (A) so all drop-down /
# non-dropdown class keys follow the same camelCase convention across sites.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "mdaTitle"          # was "MdaTitle"
$ws.Range("A8").Value = "pageTitleNewTab"   # was "pageTitlenewTab"
$ws.Range("A2").Value = "mdaTextHomePage"   # was "mdaTextHomepage"

# Match the author's final cursor position.
$ws.Range("A2").Select()
